$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("supervision")

# Insert a new row before row 6, shifting existing rows 6+ down by one.
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the Castiblanco & Wilches research project entry.
$ws.Cells.Item(6, 1).Value = "Psicología"
$ws.Cells.Item(6, 2).Value = "2022 - 2023"
$ws.Cells.Item(6, 3).Value = "Maria Camila Wilches y Johan Sebatián Castiblanco"
$ws.Cells.Item(6, 4).Value = "\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia"
$ws.Cells.Item(6, 5).Value = "Trabajo de grado: \textit{\href{https://youtu.be/FlZvukFqTcc}{El rol del género en la identificación de la sociosexualidad a partir de las voces}}"

# Row height + selection to mirror the edited workbook.
$ws.Rows.Item(6).RowHeight = 30
$ws.Range("E6").Select() | Out-Null
